$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column K: "Driver id" header with values for data rows
$ws.Range("K1").Value = "Driver id"
$ws.Range("K2").Value = 1
$ws.Range("K3").Value = 1
